$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'29.713.98"
$ws.Cells.Item(2, 5).Value = '  +3.89%  '
# Row 3
$ws.Cells.Item(3, 4).Value = "'1.912.21"
$ws.Cells.Item(3, 5).Value = '  +1.81%  '
# Row 4
$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 5).Value = '  -0.94%  '
# Row 5
$ws.Cells.Item(5, 4).Value = "'317.40"
$ws.Cells.Item(5, 5).Value = '  +0.37%  '
# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.79%  '
# Row 7
$ws.Cells.Item(7, 4).Value = "'0.5193"
$ws.Cells.Item(7, 5).Value = '  +1.65%  '
# Row 8
$ws.Cells.Item(8, 4).Value = "'0.3975"
$ws.Cells.Item(8, 5).Value = '  +1.21%  '
# Row 9
$ws.Cells.Item(9, 4).Value = "'0.08532"
$ws.Cells.Item(9, 5).Value = '  +1.58%  '
# Row 10
$ws.Cells.Item(10, 2).Value = 'Polygon'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(10, 4).Value = "'1.124"
$ws.Cells.Item(10, 5).Value = '  +0.85%  '
# Row 11
$ws.Cells.Item(11, 2).Value = 'OKB'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(11, 4).Value = "'42.59"
$ws.Cells.Item(11, 5).Value = '  +2.02%  '
# Row 12
$ws.Cells.Item(12, 4).Value = "'6.314"
$ws.Cells.Item(12, 5).Value = '  +0.54%  '
# Row 13
$ws.Cells.Item(13, 4).Value = "'1.906.06"
$ws.Cells.Item(13, 5).Value = '  +1.19%  '
# Row 14
$ws.Cells.Item(14, 4).Value = "'20.95"
$ws.Cells.Item(14, 5).Value = '  +2.25%  '
# Row 15
$ws.Cells.Item(15, 4).Value = "'7.364"
$ws.Cells.Item(15, 5).Value = '  +1.16%  '
# Row 16
$ws.Cells.Item(16, 5).Value = '  -0.80%  '
# Row 17
$ws.Cells.Item(17, 4).Value = "'93.95"
$ws.Cells.Item(17, 5).Value = '  +2.59%  '
# Row 18
$ws.Cells.Item(18, 4).Value = "'0.00001117"
$ws.Cells.Item(18, 5).Value = '  +0.93%  '
# Row 19
$ws.Cells.Item(19, 4).Value = "'0.06749"
$ws.Cells.Item(19, 5).Value = '  +0.33%  '
# Row 20
$ws.Cells.Item(20, 4).Value = "'17.99"
$ws.Cells.Item(20, 5).Value = '  +1.37%  '
# Row 21
$ws.Cells.Item(21, 5).Value = '  -0.73%  '
# Row 22
$ws.Cells.Item(22, 4).Value = "'6.044"
$ws.Cells.Item(22, 5).Value = '  +1.13%  '
# Row 23
$ws.Cells.Item(23, 4).Value = "'29.698.16"
$ws.Cells.Item(23, 5).Value = '  +3.73%  '
# Row 24
$ws.Cells.Item(24, 4).Value = "'11.22"
$ws.Cells.Item(24, 5).Value = '  +0.68%  '
# Row 25
$ws.Cells.Item(25, 4).Value = "'2.208"
$ws.Cells.Item(25, 5).Value = '  -1.73%  '
# Row 26
$ws.Cells.Item(26, 4).Value = "'2.124.17"
$ws.Cells.Item(26, 5).Value = '  +1.19%  '
# Row 27
$ws.Cells.Item(27, 4).Value = "'21.01"
$ws.Cells.Item(27, 5).Value = '  +0.71%  '
# Row 28
$ws.Cells.Item(28, 4).Value = "'159.05"
$ws.Cells.Item(28, 5).Value = '  -1.44%  '
# Row 29
$ws.Cells.Item(29, 4).Value = "'2.457"
$ws.Cells.Item(29, 5).Value = '  +3.52%  '
# Row 30
$ws.Cells.Item(30, 4).Value = "'128.35"
$ws.Cells.Item(30, 5).Value = '  +0.43%  '
# Row 31
$ws.Cells.Item(31, 4).Value = "'1.080"
$ws.Cells.Item(31, 5).Value = '  +2.40%  '
# Row 32
$ws.Cells.Item(32, 4).Value = "'0.1057"
$ws.Cells.Item(32, 5).Value = '  +0.19%  '
# Row 33
$ws.Cells.Item(33, 4).Value = "'6.226"
$ws.Cells.Item(33, 5).Value = '  +7.09%  '
# Row 34
$ws.Cells.Item(34, 4).Value = "'3.692"
$ws.Cells.Item(34, 5).Value = '  +2.14%  '
# Row 35
$ws.Cells.Item(35, 4).Value = "'0.02501"
$ws.Cells.Item(35, 5).Value = '  +1.45%  '
# Row 36
$ws.Cells.Item(36, 4).Value = "'0.06646"
$ws.Cells.Item(36, 5).Value = '  +1.47%  '
# Row 37
$ws.Cells.Item(37, 4).Value = "'9.114"
$ws.Cells.Item(37, 5).Value = '  +2.27%  '
# Row 38
$ws.Cells.Item(38, 4).Value = "'0.2209"
$ws.Cells.Item(38, 5).Value = '  +0.97%  '
# Row 39
$ws.Cells.Item(39, 4).Value = "'1.244"
$ws.Cells.Item(39, 5).Value = '  +3.82%  '
# Row 40
$ws.Cells.Item(40, 4).Value = "'5.229"
$ws.Cells.Item(40, 5).Value = '  +3.10%  '
# Row 41
$ws.Cells.Item(41, 4).Value = "'0.6567"
$ws.Cells.Item(41, 5).Value = '  +1.33%  '
# Row 42
$ws.Cells.Item(42, 4).Value = "'1.238"
$ws.Cells.Item(42, 5).Value = '  -2.29%  '
# Row 43
$ws.Cells.Item(43, 4).Value = "'11.36"
$ws.Cells.Item(43, 5).Value = '  +1.51%  '
# Row 44
$ws.Cells.Item(44, 4).Value = "'0.6143"
$ws.Cells.Item(44, 5).Value = '  +1.06%  '
# Row 45
$ws.Cells.Item(45, 4).Value = "'13.21"
$ws.Cells.Item(45, 5).Value = '  +0.95%  '
# Row 46
$ws.Cells.Item(46, 4).Value = "'3.683"
$ws.Cells.Item(46, 5).Value = '  -0.58%  '
# Row 47
$ws.Cells.Item(47, 4).Value = "'2.069"
$ws.Cells.Item(47, 5).Value = '  +1.59%  '
# Row 48
$ws.Cells.Item(48, 4).Value = "'1.240"
$ws.Cells.Item(48, 5).Value = '  +1.70%  '
# Row 49
$ws.Cells.Item(49, 4).Value = "'124.84"
$ws.Cells.Item(49, 5).Value = '  +1.86%  '
# Row 50
$ws.Cells.Item(50, 4).Value = "'1.167"
$ws.Cells.Item(50, 5).Value = '  -2.27%  '
# Row 51
$ws.Cells.Item(51, 4).Value = "'78.16"
$ws.Cells.Item(51, 5).Value = '  +1.25%  '
